$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: clear the "Total Duration:" / "0 Hours" labels -> now blank numeric cells
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = $null

# Row 4: fill in the clock-out time and the (negative) computed duration
$ws.Range("C4").Value = "11:41:38"
$ws.Range("D4").Value = "-10.36 Hours"

# Row 5 is new: another clock-in entry plus the updated running Total Duration.
# Force A5 to be stored as text first (otherwise "2026-01-23" gets parsed as a date).
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2026-01-23"
$ws.Range("B5").Value = "11:41:51"
$ws.Range("C5").Value = "Total Duration:"
$ws.Range("D5").Value = "-10.5 Hours"

# Match the formatting of the surrounding data rows (style index 2) for every
# touched cell, without introducing brand-new style/font/numFmt entries.
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)
